# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that any literal "System" token is moved to the end of its comma-separated
# list, while the relative order of all other tokens (including the
# lower-case "system" token, which is distinct from "System") is preserved.
#
# Examples:
#   "system, System, backup@backdoor.com" -> "system, backup@backdoor.com, System"
#   "System, dnasr281@gmail.com"           -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com"          -> "backup@backdoor.com, System"
#
# This engine's `-eq`/`-ceq`/`-cmatch` operators are case-insensitive, so an
# explicit, character-code based comparison is used to tell "System" (capital
# S) apart from "system" (lower-case s).

function Test-ExactEqual($a, $b) {
    if ($a.Length -ne $b.Length) { return $false }
    for ($i = 0; $i -lt $a.Length; $i++) {
        $ca = [int][char]$a.Substring($i, 1)
        $cb = [int][char]$b.Substring($i, 1)
        if ($ca -ne $cb) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$changed = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        $nonSystem = @()
        $systemCount = 0

        foreach ($p in $parts) {
            if (Test-ExactEqual $p "System") {
                $systemCount++
            } else {
                $nonSystem += $p
            }
        }

        if ($systemCount -gt 0) {
            $newParts = $nonSystem
            for ($i = 0; $i -lt $systemCount; $i++) {
                $newParts += "System"
            }
            $newVal = [string]::Join(", ", $newParts)

            if ($newVal -cne $val) {
                $cell.Value = $newVal
                $changed++
            }
        }
    }
}

Write-Host "Reordered 'System' token in $changed 'Recorded By' cell(s)."
